$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
Write-Host $ws.Range("D7").Value2
Write-Host $ws.Range("E7").Value2
